$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.848.32'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.77%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.098.45'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.49%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.41%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.14%  '

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.39%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.38'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.44%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.16'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.25%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.90%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0765'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.13%  '

# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.23%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.943'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.21%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.07'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -7.05%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.406.95'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.60%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.50'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.77%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.112.80'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.17%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.822.19'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.79%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.20'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -7.14%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.80'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.19%  '

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.24%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.48'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.95%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.46'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.61%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.66'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.84%  '

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.90%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.42'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '21.02'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.44%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.82%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.25'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.89%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.03%  '

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.54%  '

# Row 34
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.79%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.54%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.20%  '

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.28%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.26%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.27'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.89%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.02%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0221'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.05%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.90'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0961'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.25%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '96.70'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.94%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -8.26%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.89'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +16.19%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.419.32'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +11.78%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.13'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -8.28%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.40%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.07%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.293.82'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.78%  '
